# get scaffolding in place for graph loading
#
# - promote the scaffolding sheet ".EdgeCollection2" into a real, visible
#   sheet by dropping its leading-dot "private" naming convention
# - make it the active sheet/tab (was EdgeCollection1) with a fresh
#   selection at B25

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(".EdgeCollection2")
$ws.Name = "EdgeCollection2"

$ws.Activate() | Out-Null
$ws.Range("B25").Select() | Out-Null
